$wb = $excel.ActiveWorkbook

# Update the "Date" value on the Metadata sheet
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B8").Value = "2025-07-21T11:52:46+00:00"

# Update the System URI values on the Include sheets
$include0 = $wb.Worksheets.Item("Include #0")
$include0.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R262-CategorieSocioProfessionnelle"

$include1 = $wb.Worksheets.Item("Include #1")
$include1.Range("B4").Value = "https://interop.esante.gouv.fr/terminologies/CodeSystem-TRE-R292-INSEECategorieSocioProfessionnelleAgrNiv1"
